$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the three new log rows under the existing "Time / Player / Coin" header.
$ws.Range("A2").Value = "Sat Dec 16 23_40_40 2023"
$ws.Range("B2").Value = "loc"
$ws.Range("C2").Value = -10

$ws.Range("A3").Value = "Sat Dec 16 23_42_51 2023"
$ws.Range("B3").Value = "test"
$ws.Range("C3").Value = 50

$ws.Range("A4").Value = "Sat Dec 16 23_49_37 2023"
$ws.Range("B4").Value = "test"
$ws.Range("C4").Value = -10

# Mirror the author's final selection (A2:H3, anchored so H3 is reachable).
$ws.Range("A2:H3").Select() | Out-Null
